# Auto-generated Excel COM-interop script to update the cryptos list
# with refreshed Price (column D) and Volume(1h) (column E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.687.69"
$ws.Range("E2").Value = "  -3.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.550.61"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.54"
$ws.Range("E5").Value = "  -3.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.63"
$ws.Range("E6").Value = "  -7.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.555"
$ws.Range("E8").Value = "  -5.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.554.20"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("E10").Value = "  -6.92%  "
$ws.Range("E11").Value = "  -4.70%  "
$ws.Range("E12").Value = "  -4.88%  "
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.004.50"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.677.46"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.47"
$ws.Range("E16").Value = "  -5.35%  "
$ws.Range("E17").Value = "  -5.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.562.05"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.52"
$ws.Range("E19").Value = "  -4.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "332.70"
$ws.Range("E20").Value = "  -6.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").Value = "  -5.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  -4.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.79"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.405"
$ws.Range("E25").Value = "  -4.98%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  -3.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0777"
$ws.Range("E28").Value = "  -8.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.85"
$ws.Range("E29").Value = "  -7.35%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "149.75"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.51"
$ws.Range("E32").Value = "  -4.88%  "
$ws.Range("E33").Value = "  -4.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  -8.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.87"
$ws.Range("E35").Value = "  -6.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.896"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.10"
$ws.Range("E37").Value = "  -8.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.83"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.820"
$ws.Range("E39").Value = "  -8.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "286.47"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  -7.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.50"
$ws.Range("E42").Value = "  -7.40%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0983"
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0530"
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.50"
$ws.Range("E48").Value = "  -5.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"
$ws.Range("E49").Value = "  -4.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.51"
$ws.Range("E50").Value = "  -8.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.890.57"
$ws.Range("E51").Value = "  -4.16%  "
